$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 / row 2 (D column)
$ws.Range("D1").Value = "DurationTime"
$ws.Range("D2").Value = "완료시간(이동시간)"

# Update the "Left" row (row 4) values
$ws.Range("A4").Value = "Left"
$ws.Range("B4").Value = -10
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 2

# Update the "Right" row (row 5) values
$ws.Range("A5").Value = "Right"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 2

# Remove the now-unused rows 6-10 (LeftSlow/Right/LeftFast/RightSlow/Pause)
$ws.Rows("6:10").Delete()

# Restore the previously-saved selection
$ws.Range("B27").Select()
